$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 21 (NiLatticeSlosh) ---
# RoiName (D21): Bec -> NiLattice
$ws.Range("D21").Value = "NiLattice"
# AnalysisMethod (G21): remove "Tof;" segment
$ws.Range("G21").Value = "DensityFit;AtomNumber;CenterFit"
# CenterFitMethod (N21): ParabolicFit1D -> SineFit1D
$ws.Range("N21").Value = "SineFit1D"

# --- Add new row 22 (NiLattice) ---
$ws.Range("A22").Value = "NiLattice"
$ws.Range("B22").Value = "A slosh experiment at the non-interacting lattice stage."
$ws.Range("C22").Value = "TOP"
$ws.Range("D22").Value = "NiLattice"
$ws.Range("E22").Value = 4
$ws.Range("F22").Value = "RunIndex"
$ws.Range("G22").Value = "DensityFit;AtomNumber"
$ws.Range("H22").Value = "LSR"
$ws.Range("I22").Value = "HF"
$ws.Range("J22").Value = "StrongLight"
$ws.Range("K22").Value = 8
$ws.Range("L22").Value = "BosonicGaussianFit1D"
$ws.Range("M22").Value = 1
$ws.Range("N22").Value = "ParabolicFit1D"

# --- Add new row 23 (NiBec) ---
$ws.Range("A23").Value = "NiBec"
$ws.Range("B23").Value = "A TOF experiment at the non-interacting BEC stage."
$ws.Range("C23").Value = "TOP"
$ws.Range("D23").Value = "Bec"
$ws.Range("E23").Value = 4
$ws.Range("F23").Value = "RunIndex"
$ws.Range("G23").Value = "CenterFit;AtomNumber;DensityFit"
$ws.Range("H23").Value = "LSR"
$ws.Range("I23").Value = "HF"
$ws.Range("J23").Value = "StrongLight"
$ws.Range("K23").Value = 8
$ws.Range("L23").Value = "BosonicGaussianFit1D"
$ws.Range("M23").Value = 1
$ws.Range("N23").Value = "ParabolicFit1D"
